# Bold the text in the last row of every table in the deck.
#
# Affected slides (by position in $p.Slides, which matches the
# ppt/slides/slideN.xml numbering 1:1 for this deck):
#   4  - "AWS Services Deployed" row
#   9  - "Wave 3" row
#   11 - "Monitoring" row
#   15 - "Productivity Gains" row
#   21 - "Change Management Resistance" row
#   26 - "AI/ML Service Integration" row

$p = $ppt.ActivePresentation

$slideIndexes = @(4, 9, 11, 15, 21, 26)

foreach ($slideIdx in $slideIndexes) {
    $s = $p.Slides.Item($slideIdx)

    for ($shpIdx = 1; $shpIdx -le $s.Shapes.Count; $shpIdx++) {
        $shp = $s.Shapes.Item($shpIdx)

        if ($shp.HasTable) {
            $tbl = $shp.Table
            $lastRow = $tbl.Rows.Count

            for ($col = 1; $col -le $tbl.Columns.Count; $col++) {
                $cell = $tbl.Cell($lastRow, $col)
                $cell.Shape.TextFrame.TextRange.Font.Bold = 1
            }
        }
    }
}
